$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number formats that map back onto the workbook's existing styles:
#  - "yyyy\-mm\-dd"  -> custom numFmtId 164 (style index 1 in the original file)
#  - "m/d/yy"        -> builtin numFmtId 14  (style index 2 in the original file)
$fmtCustom  = "yyyy\-mm\-dd"
$fmtBuiltin = "m/d/yy"

# --- Row 4: second (first split) term for HAN-TOPMODEL-002, comprising 2012-2013 ---
$ws.Range("A4").Value() = "HAN-TOPMODEL-002"
$ws.Range("B4").Value() = "4060-0090-03"

$ws.Range("C4").Value() = 40544
$ws.Range("C4").NumberFormat() = $fmtBuiltin

$ws.Range("D4").Value() = 40908
$ws.Range("D4").NumberFormat() = $fmtCustom

$ws.Range("E4").Value() = 2100

$ws.Range("F4").Value() = 40909
$ws.Range("F4").NumberFormat() = $fmtBuiltin

$ws.Range("G4").Value() = 41639
$ws.Range("G4").NumberFormat() = $fmtBuiltin

$ws.Range("H4").Value() = 2200
$ws.Range("I4").Value() = 3.6
$ws.Range("J4").Value() = 2012

# --- Row 5: third (second split) term for HAN-TOPMODEL-002, comprising 2014 ---
$ws.Range("A5").Value() = "HAN-TOPMODEL-002"
$ws.Range("B5").Value() = "4060-0090-03"

$ws.Range("C5").Value() = 41275
$ws.Range("C5").NumberFormat() = $fmtBuiltin

$ws.Range("D5").Value() = 41639
$ws.Range("D5").NumberFormat() = $fmtBuiltin

$ws.Range("E5").Value() = 2300

$ws.Range("F5").Value() = 41640
$ws.Range("F5").NumberFormat() = $fmtBuiltin

$ws.Range("G5").Value() = 42004
$ws.Range("G5").NumberFormat() = $fmtBuiltin

$ws.Range("H5").Value() = 2300
$ws.Range("J5").Value() = 2014

# Match the saved selection state left behind in the authored workbook.
$ws.Range("K8").Select() | Out-Null
